$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append row 4 (new inbound e-mail, not auto-answered) ---
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Range("A4").Value = "Interne taak"
$wsLogs.Range("B4").Value = "kwaliteit@testbedrijf123.nl"
$wsLogs.Range("C4").Value = "Leg dit even neer bij Koen."
$wsLogs.Range("D4").Value = "Onbekend"
# E4 ("Antwoord") intentionally left blank - no automatic reply was sent.
$wsLogs.Range("F4").Value = "2025-08-18 20:20:28"
$wsLogs.Range("G4").Value = "Nee"
$wsLogs.Range("H4").Value = "Ja"
$wsLogs.Range("I4").Value = "Nee"
$wsLogs.Range("J4").Value = "Nee"

# Extend the existing conditional-formatting rules to cover the new row
# without touching their dxfs/priorities.
$wsLogs.Range("D2:D3").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("D2:D4"))
$wsLogs.Range("G2:G3").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("G2:G4"))
$wsLogs.Range("H2:H3").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("H2:H4"))
$wsLogs.Range("I2:I3").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("I2:I4"))
$wsLogs.Range("J2:J3").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("J2:J4"))

# --- "Dashboard" sheet: add the summary row for the new "Onbekend" category ---
$wsDash = $wb.Worksheets.Item("Dashboard")

$wsDash.Range("A3").Value = "Onbekend"
$wsDash.Range("B3").Value = 1

# --- Chart: grow the category/value series references to include row 3 ---
$chart = $wsDash.ChartObjects().Item(1).Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$3,'Dashboard'!`$B`$2:`$B`$3,1)"
